# Updates in ELN Setup
# 1) Fix typo in first worksheet's name
# 2) Fix typo in the shared "hasStoichiometriCoefficient" label used on
#    sheet 1 ("Substances and Parameters"), cell A5.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Substances and Parameters"
$ws1.Range("A5").Value = "hasStoichiometricCoefficient"

$ws2 = $wb.Worksheets.Item(2)

# Restore / update the active-cell selections on sheet 1 and sheet 2
[void]$ws1.Activate()
[void]$ws1.Range("C14").Select()
[void]$ws2.Activate()
[void]$ws2.Range("B5").Select()
[void]$ws1.Activate()
